# "se crea punto 42" - add the new quarterly (2nd trimester 2023) report row
# to the "Reporte de Formatos" sheet: shift the reporting period from
# Q1 2023 (01/01-31/03) to Q2 2023 (01/04-30/06), push out the
# validation/update dates, and expand the justification note in column M.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# --- Period covered by the report: 01/04/2023 - 30/06/2023 ---
$ws.Range("B8").Value = (Get-Date -Year 2023 -Month 4 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C8").Value = (Get-Date -Year 2023 -Month 6 -Day 30 -Hour 0 -Minute 0 -Second 0).Date

# --- Validation / update dates: 10/07/2023 ---
$ws.Range("K8").Value = (Get-Date -Year 2023 -Month 7 -Day 10 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("L8").Value = (Get-Date -Year 2023 -Month 7 -Day 10 -Hour 0 -Minute 0 -Second 0).Date

# --- Expanded note explaining the lack of a Consejo Consultivo ---
$ws.Range("M8").Value = "La Universidad Politécnica de Pachuca, no tiene Consejo Consultivo por que se encuentra en proceso de validación el nuevo Decreto de Creación."
$ws.Range("M8").HorizontalAlignment = -4130
$ws.Range("M8").WrapText = $true
$ws.Columns.Item(13).ColumnWidth = 75.57

# --- Row heights grow to fit the longer wrapped note / header ---
$ws.Rows.Item(3).RowHeight = 63
$ws.Rows.Item(8).RowHeight = 30
